$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shift the existing rows 39-44 down to rows 40-45, making room at row 39
#    for the new LONP1 gene entry (columns A-E).
$ws.Range("A39:E44").Copy()
$ws.Range("A40:E45").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# The paste above loses the style on the very last row of the pasted
# range (A45), since it previously lay just beyond the sheet's used range.
# Re-apply A44's style (bold/centered/bordered, same as every other cell
# in column A) onto A45 to restore it, without touching its value.
$ws.Range("A44").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2. Populate the now-empty row 39 with the LONP1 gene entry.
#    Copy the D38 cell (a text "1") into D39 first so the geneConfidence
#    value stays a text "1" instead of being auto-coerced to a number.
$ws.Range("D38").Copy()
$ws.Range("D39").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Cells.Item(39, 2).Value = "LONP1"
$ws.Cells.Item(39, 3).Value = "lon peptidase 1, mitochondrial"
$ws.Cells.Item(39, 5).Value = "Congenital diaphragmatic hernia"

# 3. Add the new "time_taken" header in F1, copying the style of the
#    neighboring header cell E1 so it matches (bold, bordered, centered).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(1, 6).Value = "time_taken"

# 4. Populate F2:F45 with the time_taken timestamps recorded for each gene.
$timestamps = @(
  "2021-10-05 10:50:35.521724",
  "2021-10-05 10:50:35.521734",
  "2021-10-05 10:50:35.521738",
  "2021-10-05 10:50:35.521740",
  "2021-10-05 10:50:35.521743",
  "2021-10-05 10:50:35.521746",
  "2021-10-05 10:50:35.521749",
  "2021-10-05 10:50:35.521751",
  "2021-10-05 10:50:35.521754",
  "2021-10-05 10:50:35.521757",
  "2021-10-05 10:50:35.521759",
  "2021-10-05 10:50:35.521762",
  "2021-10-05 10:50:35.521764",
  "2021-10-05 10:50:35.521767",
  "2021-10-05 10:50:35.521769",
  "2021-10-05 10:50:35.521772",
  "2021-10-05 10:50:35.521774",
  "2021-10-05 10:50:35.521777",
  "2021-10-05 10:50:35.521780",
  "2021-10-05 10:50:35.521782",
  "2021-10-05 10:50:35.521785",
  "2021-10-05 10:50:35.521787",
  "2021-10-05 10:50:35.521790",
  "2021-10-05 10:50:35.521792",
  "2021-10-05 10:50:35.521795",
  "2021-10-05 10:50:35.521798",
  "2021-10-05 10:50:35.521801",
  "2021-10-05 10:50:35.521803",
  "2021-10-05 10:50:35.521806",
  "2021-10-05 10:50:35.521808",
  "2021-10-05 10:50:35.521811",
  "2021-10-05 10:50:35.521814",
  "2021-10-05 10:50:35.521817",
  "2021-10-05 10:50:35.521820",
  "2021-10-05 10:50:35.521822",
  "2021-10-05 10:50:35.521825",
  "2021-10-05 10:50:35.521827",
  "2021-10-05 10:50:35.521830",
  "2021-10-05 10:50:35.521832",
  "2021-10-05 10:50:35.521835",
  "2021-10-05 10:50:35.521838",
  "2021-10-05 10:50:35.521841",
  "2021-10-05 10:50:35.521843",
  "2021-10-05 10:50:35.521846"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

Write-Host "Edit complete"
